$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added for Membrillo (Vega Modelo de
# Temuco). It belongs right after the header/first block, at row 141,
# pushing the existing rows 141-167 down to 142-168.
$ws.Rows(141).Insert()

$ws.Range("A141").Value = 10
$ws.Range("B141").Value = "Vega Modelo de Temuco"
$ws.Range("C141").Value = "La Araucanía"
$ws.Range("D141").Value = 44694
$ws.Range("E141").Value = 9
$ws.Range("F141").Value = "Fruta"
$ws.Range("G141").Value = 100104
$ws.Range("H141").Value = "Frutos de pepita"
$ws.Range("I141").Value = 100104003
$ws.Range("J141").Value = "Membrillo"
$ws.Range("K141").Value = "Champion"
$ws.Range("L141").Value = "Primera"
$ws.Range("M141").Value = 95
$ws.Range("N141").Value = 13000
$ws.Range("O141").Value = 13000
$ws.Range("P141").Value = 13000
$ws.Range("Q141").Value = "$/bandeja 18 kilos granel"
$ws.Range("R141").Value = "Región de O'Higgins"
$ws.Range("S141").Value = 722
$ws.Range("T141").Value = 18
